$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: new working-hours entry
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats (reuse the existing date style)
$ws.Range("A15").Value = 45428
$ws.Range("B15").Value = 5
$ws.Range("D15").Value = "added prediction metrics"

# Row 16: TODO / links note
$ws.Range("F16").Value = "Look into other and all metrics, the potential of RF's again"

# Update selection to match the saved workbook state
$ws.Range("F16").Select()
